# bom_panel.xlsx — "finished schematic, started panel BOM"
# Populate the panel bill-of-materials grid (Item / MPN / DPN / Quantity /
# Description / Notes) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "MPN"
$ws.Range("C1").Value = "DPN"
$ws.Range("D1").Value = "Quantity"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Notes"

# --- AUX Connector section --------------------------------------------
$ws.Range("A2").Value = "AUX Connector"

$ws.Range("A3").Value = "Power distribution wire"
$ws.Range("B3").Value = "55A0121-12-9/96CS2275"
$ws.Range("C3").Value = "55A0121-12-9/96CS2275-DS-ND"
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "2 (1 Pair Twisted) Conductor Multi-Conductor Cable  12 AWG  Enter Number of Feet in Order Quantity"
$ws.Range("F3").Value = "qt in feet"

# --- Standoffs section --------------------------------------------------
$ws.Range("A4").Value = "Standoffs"

$ws.Range("A5").Value = "Input wire"
$ws.Range("B5").Value = "E2102S.41.02"
$ws.Range("C5").Value = "CE2102W-25-ND"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "2 Conductor Multi-Conductor Cable White 22 AWG Foil 25.00' (7.62m)"
$ws.Range("F5").Value = "7/30 strands"

$ws.Range("A6").Value = "heat shrink (Audio input)"
$ws.Range("B6").Value = "V2-1.5-0-SP-SM"
$ws.Range("C6").Value = "V2-1.5-R5-ND"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 'Heat Shrink Tubing, Flexible 0.083" (2.11mm) 2 to 1 Black 0.042'' (12.70mm, 0.50")'

$ws.Range("A7").Value = "Potentiometer"
$ws.Range("B7").Value = "EVU-F3MFL3D14"
$ws.Range("C7").Value = "P3F6103-ND"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = "10k Ohm 1 Gang Logarithmic Panel Mount Potentiometer Detent 1 Kierros  0.05W, 1/20W PC Pins, Board Locks"

$ws.Range("A8").Value = "Power Receptacle"
$ws.Range("B8").Value = "KM01.1205.11"
$ws.Range("C8").Value = "486-2280-ND"
$ws.Range("E8").Value = "Power Entry Connector Receptacle, Male Blades - Module IEC 320-C14 Panel Mount, Snap-In"
$ws.Range("F8").Value = "4.8 mm Quick Connect"

# DPN column on this row is shown in a slightly larger Arial font
$ws.Range("C8").Font.Size = 13
$ws.Range("C8").Font.Name = "Arial"

$ws.Range("A9").Value = "heat shrink (power input)"

$ws.Range("A10").Value = "Blade Connector (power input)"

# --- Column sizing (auto-fit to content, like Excel's "AutoFit Column
# Width") ----------------------------------------------------------------
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).ColumnWidth = 18
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()

# --- Selection, matching the saved workbook view -----------------------
$ws.Range("B29").Select()
